# "atualizei dados da add"
# - Update the existing June 13th total_venda value (row 11).
# - Insert a new daily record for June 16th right after it, pushing every
#   subsequent row (old rows 12..71) down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B11 (Dia 13, Mes 6) total_venda: 6825.70.. -> 8459.48
$ws.Range("B11").Value = 8459.48

# Insert a new blank row at position 12 (shifts old rows 12-71 -> 13-72)
$ws.Rows.Item(12).Insert()

# Fill the newly inserted row 12 with the new June 16th record
$ws.Range("A12").Value = 16
$ws.Range("B12").Value = 23567.89
$ws.Range("C12").Value = 6
$ws.Range("D12").Value = 2025
$ws.Range("E12").Value = "06/2025"
